# Update the "Final ML Models" sheet with revised LSTM multivariate model results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final ML Models")

$ws.Range("B8").Value = "lSTM Multivariate with 1 lag"
$ws.Range("C8").Value = "RMSE= 0.095"
$ws.Range("D8").Value = "RMSE= 0.054"
$ws.Range("E8").Value = "RMSE= 0.094"

$ws.Range("C9").Select()
